$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 3 new "Spiral" rows right after the row that will hold the
# Gaussian-Quadrature data (row 10). This shifts the old rows 11-16 down to rows 14-19,
# carrying their formatting (including the bold/bordered style used in column A) with them.
$ws.Rows("11:13").Insert()

# The 3 freshly inserted rows don't pick up the same column-A style as the rest of the
# table, so copy it over from a neighboring, already-styled cell.
$ws.Range("A10").Copy() | Out-Null
$ws.Range("A11:A13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Row 10: Gaussian-Quadrature (previously the last row of the block) ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.016397467914081
$ws.Cells.Item(10, 4).Value = 1.164293674555247
$ws.Cells.Item(10, 5).Value = 0.9387573940798133
$ws.Cells.Item(10, 6).Value = 1.016397467914081
$ws.Cells.Item(10, 7).Value = 1.057764820783397
$ws.Cells.Item(10, 8).Value = 0.9383858076508292
$ws.Cells.Item(10, 9).Value = 0.9584598127781756
$ws.Cells.Item(10, 10).Value = 1.164293674555247
$ws.Cells.Item(10, 11).Value = 1.05152553431753
$ws.Cells.Item(10, 12).Value = 1.033961501115806
$ws.Cells.Item(10, 13).Value = 1.012343162960257

# --- Row 11: Spiral-90deg-10rot-5space (new) ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 0.8651284232307052
$ws.Cells.Item(11, 4).Value = 0.9491879909589457
$ws.Cells.Item(11, 5).Value = 1.097094467977835
$ws.Cells.Item(11, 6).Value = 0.8651284232307052
$ws.Cells.Item(11, 7).Value = 0.8564506183700774
$ws.Cells.Item(11, 8).Value = 1.420787206736747
$ws.Cells.Item(11, 9).Value = 1.031116774713991
$ws.Cells.Item(11, 10).Value = 0.9491879909589457
$ws.Cells.Item(11, 11).Value = 1.02314122946839
$ws.Cells.Item(11, 12).Value = 0.9441348263495477
$ws.Cells.Item(11, 13).Value = 1.036627580331384

# --- Row 12: Spiral-90deg-15rot-5space (new) ---
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 0.8648546335092328
$ws.Cells.Item(12, 4).Value = 0.952154130189748
$ws.Cells.Item(12, 5).Value = 1.096198270893054
$ws.Cells.Item(12, 6).Value = 0.8648546335092328
$ws.Cells.Item(12, 7).Value = 0.8580951606762712
$ws.Cells.Item(12, 8).Value = 1.418456821763833
$ws.Cells.Item(12, 9).Value = 1.030366014182722
$ws.Cells.Item(12, 10).Value = 0.952154130189748
$ws.Cells.Item(12, 11).Value = 1.024176200541401
$ws.Cells.Item(12, 12).Value = 0.9445154170253169
$ws.Cells.Item(12, 13).Value = 1.036687505202477

# --- Row 13: Spiral-90deg-10rot-3space (new) ---
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 0.8647370851650312
$ws.Cells.Item(13, 4).Value = 0.950433596945937
$ws.Cells.Item(13, 5).Value = 1.096875468137782
$ws.Cells.Item(13, 6).Value = 0.8647370851650312
$ws.Cells.Item(13, 7).Value = 0.8568942790295715
$ws.Cells.Item(13, 8).Value = 1.420693219209383
$ws.Cells.Item(13, 9).Value = 1.030851293857961
$ws.Cells.Item(13, 10).Value = 0.950433596945937
$ws.Cells.Item(13, 11).Value = 1.02365453254186
$ws.Cells.Item(13, 12).Value = 0.9441958088534455
$ws.Cells.Item(13, 13).Value = 1.036747490390944

# --- Row 14: NoRotation-tilt60deg (values that used to live in row 10) ---
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 0.7060799999999997
$ws.Cells.Item(14, 4).Value = 1.338312000000001
$ws.Cells.Item(14, 5).Value = 1.000016
$ws.Cells.Item(14, 6).Value = 0.7060799999999997
$ws.Cells.Item(14, 7).Value = 1.178648000000001
$ws.Cells.Item(14, 8).Value = 0.9008359999999992
$ws.Cells.Item(14, 9).Value = 0.9190360000000001
$ws.Cells.Item(14, 10).Value = 1.338312000000001
$ws.Cells.Item(14, 11).Value = 1.169164
$ws.Cells.Item(14, 12).Value = 0.937622
$ws.Cells.Item(14, 13).Value = 1.007154666666667

# --- Row 15: Rotation-NoTilt (values that used to live in row 11) ---
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 0.7
$ws.Cells.Item(15, 4).Value = 1.293887499999999
$ws.Cells.Item(15, 5).Value = 0.9901374999999997
$ws.Cells.Item(15, 6).Value = 0.7
$ws.Cells.Item(15, 7).Value = 1.260825000000001
$ws.Cells.Item(15, 8).Value = 0.6504250000000005
$ws.Cells.Item(15, 9).Value = 0.91
$ws.Cells.Item(15, 10).Value = 1.293887499999999
$ws.Cells.Item(15, 11).Value = 1.142012499999999
$ws.Cells.Item(15, 12).Value = 0.9210062499999996
$ws.Cells.Item(15, 13).Value = 0.9675458333333332

# --- Row 16: Rotation-60detTilt (values that used to live in row 12) ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 0.8238491901951993
$ws.Cells.Item(16, 4).Value = 1.1706334416896
$ws.Cells.Item(16, 5).Value = 0.990390712729598
$ws.Cells.Item(16, 6).Value = 0.8238491901951993
$ws.Cells.Item(16, 7).Value = 1.149711720652797
$ws.Cells.Item(16, 8).Value = 0.7965953888256005
$ws.Cells.Item(16, 9).Value = 0.9451911307263993
$ws.Cells.Item(16, 10).Value = 1.1706334416896
$ws.Cells.Item(16, 11).Value = 1.080512077209599
$ws.Cells.Item(16, 12).Value = 0.9521806337023991
$ws.Cells.Item(16, 13).Value = 0.9793952641365323

# --- Row 17: HexGrid-90degTilt5degRes (values that used to live in row 13, recomputed) ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9949362300993597
$ws.Cells.Item(17, 4).Value = 0.9950778694753355
$ws.Cells.Item(17, 5).Value = 0.9953632732079455
$ws.Cells.Item(17, 6).Value = 0.9949362300993597
$ws.Cells.Item(17, 7).Value = 0.9961725939892997
$ws.Cells.Item(17, 8).Value = 0.9959845843276111
$ws.Cells.Item(17, 9).Value = 0.9949283083447468
$ws.Cells.Item(17, 10).Value = 0.9950778694753355
$ws.Cells.Item(17, 11).Value = 0.9952205713416404
$ws.Cells.Item(17, 12).Value = 0.9950784007205
$ws.Cells.Item(17, 13).Value = 0.9954104765740497

# --- Row 18: HexGrid-90degTilt22p5degRes (values that used to live in row 14) ---
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 1.019517165461276
$ws.Cells.Item(18, 4).Value = 0.97273476138757
$ws.Cells.Item(18, 5).Value = 0.9921710918358011
$ws.Cells.Item(18, 6).Value = 1.019517165461276
$ws.Cells.Item(18, 7).Value = 0.9854260773953819
$ws.Cells.Item(18, 8).Value = 0.988653839220813
$ws.Cells.Item(18, 9).Value = 0.9986721942179078
$ws.Cells.Item(18, 10).Value = 0.97273476138757
$ws.Cells.Item(18, 11).Value = 0.9824529266116855
$ws.Cells.Item(18, 12).Value = 1.000985046036481
$ws.Cells.Item(18, 13).Value = 0.9928625215864583

# --- Row 19: HexGrid-60degTilt5degRes (values that used to live in row 15) ---
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 1.00058675045867
$ws.Cells.Item(19, 4).Value = 0.9288437873332713
$ws.Cells.Item(19, 5).Value = 1.011466811523736
$ws.Cells.Item(19, 6).Value = 1.00058675045867
$ws.Cells.Item(19, 7).Value = 0.9666629710372961
$ws.Cells.Item(19, 8).Value = 1.025298675524156
$ws.Cells.Item(19, 9).Value = 1.010030762875268
$ws.Cells.Item(19, 10).Value = 0.9288437873332713
$ws.Cells.Item(19, 11).Value = 0.9701552994285038
$ws.Cells.Item(19, 12).Value = 0.9853710249435869
$ws.Cells.Item(19, 13).Value = 0.9904816264587328

Write-Output "Applied averaged-intensities update for spiral schemes"
